{"js": "// Update the date line and the twenty-five multiplication exercises\n// (three-digit number x one-digit number) to the new day's values.\nconst replacements = [\n  [\"2025-07-20 Sunday\", \"2025-07-21 Monday\"],\n  [\"101\u00d73=\", \"409\u00d78=\"],\n  [\"138\u00d73=\", \"325\u00d72=\"],\n  [\"288\u00d77=\", \"635\u00d77=\"],\n  [\"170\u00d75=\", \"717\u00d77=\"],\n  [\"321\u00d79=\", \"514\u00d74=\"],\n  [\"587\u00d77=\", \"182\u00d74=\"],\n  [\"169\u00d72=\", \"616\u00d79=\"],\n  [\"464\u00d78=\", \"709\u00d74=\"],\n  [\"307\u00d75=\", \"481\u00d77=\"],\n  [\"757\u00d77=\", \"513\u00d77=\"],\n  [\"452\u00d72=\", \"193\u00d79=\"],\n  [\"384\u00d72=\", \"257\u00d74=\"],\n  [\"131\u00d73=\", \"234\u00d78=\"],\n  [\"460\u00d72=\", \"773\u00d79=\"],\n  [\"526\u00d74=\", \"352\u00d78=\"],\n  [\"545\u00d74=\", \"234\u00d73=\"],\n  [\"305\u00d77=\", \"189\u00d73=\"],\n  [\"173\u00d72=\", \"510\u00d74=\"],\n  [\"134\u00d72=\", \"323\u00d76=\"],\n  [\"255\u00d79=\", \"956\u00d77=\"],\n  [\"697\u00d72=\", \"722\u00d72=\"],\n  [\"650\u00d78=\", \"304\u00d73=\"],\n  [\"984\u00d74=\", \"606\u00d76=\"],\n  [\"342\u00d75=\", \"394\u00d72=\"],\n  [\"378\u00d74=\", \"670\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five multiplication exercises\n# (three-digit number x one-digit number) to the new day's values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-20 Sunday\", \"2025-07-21 Monday\"),\n    @(\"101\u00d73=\", \"409\u00d78=\"),\n    @(\"138\u00d73=\", \"325\u00d72=\"),\n    @(\"288\u00d77=\", \"635\u00d77=\"),\n    @(\"170\u00d75=\", \"717\u00d77=\"),\n    @(\"321\u00d79=\", \"514\u00d74=\"),\n    @(\"587\u00d77=\", \"182\u00d74=\"),\n    @(\"169\u00d72=\", \"616\u00d79=\"),\n    @(\"464\u00d78=\", \"709\u00d74=\"),\n    @(\"307\u00d75=\", \"481\u00d77=\"),\n    @(\"757\u00d77=\", \"513\u00d77=\"),\n    @(\"452\u00d72=\", \"193\u00d79=\"),\n    @(\"384\u00d72=\", \"257\u00d74=\"),\n    @(\"131\u00d73=\", \"234\u00d78=\"),\n    @(\"460\u00d72=\", \"773\u00d79=\"),\n    @(\"526\u00d74=\", \"352\u00d78=\"),\n    @(\"545\u00d74=\", \"234\u00d73=\"),\n    @(\"305\u00d77=\", \"189\u00d73=\"),\n    @(\"173\u00d72=\", \"510\u00d74=\"),\n    @(\"134\u00d72=\", \"323\u00d76=\"),\n    @(\"255\u00d79=\", \"956\u00d77=\"),\n    @(\"697\u00d72=\", \"722\u00d72=\"),\n    @(\"650\u00d78=\", \"304\u00d73=\"),\n    @(\"984\u00d74=\", \"606\u00d76=\"),\n    @(\"342\u00d75=\", \"394\u00d72=\"),\n    @(\"378\u00d74=\", \"670\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $r = $d.Content\n    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
